$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (data row) updates per the diff.
# Numeric-looking text values are prefixed with a leading apostrophe so
# Excel keeps storing them as text (shared strings) instead of silently
# converting them to numbers.
$ws.Range("A2").Value = "24/02/2020"
$ws.Range("B2").Value = "02:02:20"
$ws.Range("C2").Value = "'7.0"
$ws.Range("E2").Value = "nose"
$ws.Range("F2").Value = "'9.0"
$ws.Range("G2").Value = "'10.0"
$ws.Range("I2").Value = "SAN_JOSE, sirveee"
